# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 52 (the top of the
# "Mandarina" history block), pushing the existing rows 52-68 down by
# one row (they become rows 53-69). The sheet has no formulas, tables
# or autofilters below the header row, so a plain row insert followed
# by writing the new row's values reproduces the target state exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52..68 down to 53..69, leaving row 52 free for the new entry.
$ws.Rows.Item(52).Insert()

# Populate the newly freed row 52 with the new weekly record.
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").Value = 44468
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100102
$ws.Range("H52").Value = "Cítricos"
$ws.Range("I52").Value = 100102004
$ws.Range("J52").Value = "Mandarina"
$ws.Range("K52").Value = "Murcott"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 300
$ws.Range("N52").Value = 11000
$ws.Range("O52").Value = 12000
$ws.Range("P52").Value = 11500
$ws.Range("Q52").Value = "$/caja 20 kilos"
$ws.Range("R52").Value = "Región de Coquimbo"
$ws.Range("S52").Value = 575
$ws.Range("T52").Value = 20
